$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 82
$ws.Range("H82").Value = 1346.9333
$ws.Range("I82").Value = 1086
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 3258
$ws.Range("L82").Value = 15000
$ws.Range("M82").Value = -2852
$ws.Range("N82").Value = -15812
# row 85
$ws.Range("H85").Value = 1346.9333
$ws.Range("I85").Value = 1086
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 3258
$ws.Range("L85").Value = 15000
$ws.Range("M85").Value = -1854
$ws.Range("N85").Value = -17808
# row 88
$ws.Range("H88").Value = 6068.4165
$ws.Range("I88").Value = 5624
$ws.Range("J88").Value = 6290.625
$ws.Range("K88").Value = 5624
$ws.Range("L88").Value = 6290.625
$ws.Range("M88").Value = -5218
$ws.Range("N88").Value = -7102.625
# row 91
$ws.Range("H91").Value = 6068.4165
$ws.Range("I91").Value = 5624
$ws.Range("J91").Value = 6290.625
$ws.Range("K91").Value = 5624
$ws.Range("L91").Value = 6290.625
$ws.Range("M91").Value = -4220
$ws.Range("N91").Value = -9098.625
# row 132
$ws.Range("H132").Value = 5379
$ws.Range("I132").Value = 2446.4644
$ws.Range("J132").Value = 21801.2
$ws.Range("K132").Value = 7339.3932
$ws.Range("L132").Value = 65403.60000000001
$ws.Range("M132").Value = -4809.3932
$ws.Range("N132").Value = -70463.60000000001
# row 137
$ws.Range("H137").Value = 2138193.5
$ws.Range("I137").Value = 1483.8948
$ws.Range("J137").Value = 7937833.5
$ws.Range("K137").Value = 4451.6844
$ws.Range("L137").Value = 23813500.5
$ws.Range("M137").Value = -1901.6844
$ws.Range("N137").Value = -23818600.5

$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("H61").Value = 1395.8654
$ws.Range("I61").Value = 1253.4615
$ws.Range("K61").Value = 1253.4615
$ws.Range("M61").Value = -1041.4615
# row 63
$ws.Range("H63").Value = 2924.2856
$ws.Range("I63").Value = 2434.4443
$ws.Range("J63").Value = 3291.6667
$ws.Range("K63").Value = 2434.4443
$ws.Range("L63").Value = 3291.6667
$ws.Range("M63").Value = -1748.4443
$ws.Range("N63").Value = -4663.6667
# row 66
$ws.Range("H66").Value = 2924.2856
$ws.Range("I66").Value = 2434.4443
$ws.Range("J66").Value = 3291.6667
$ws.Range("K66").Value = 12172.2215
$ws.Range("L66").Value = 16458.3335
$ws.Range("M66").Value = -8740.2215
$ws.Range("N66").Value = -23322.3335
# row 74
$ws.Range("H74").Value = 16478.338
$ws.Range("I74").Value = 18962.072
$ws.Range("J74").Value = 2817.8
$ws.Range("K74").Value = 18962.072
$ws.Range("L74").Value = 2817.8
$ws.Range("M74").Value = -18088.072
$ws.Range("N74").Value = -4565.8
# row 77
$ws.Range("H77").Value = 16478.338
$ws.Range("I77").Value = 18962.072
$ws.Range("J77").Value = 2817.8
$ws.Range("K77").Value = 94810.36
$ws.Range("L77").Value = 14089
$ws.Range("M77").Value = -90442.36
$ws.Range("N77").Value = -22825
# row 88
$ws.Range("H88").Value = 2670
$ws.Range("I88").Value = 2980
$ws.Range("J88").Value = 2566.6667
$ws.Range("K88").Value = 2980
$ws.Range("L88").Value = 2566.6667
$ws.Range("M88").Value = -2574
$ws.Range("N88").Value = -3378.6667
# row 91
$ws.Range("H91").Value = 2670
$ws.Range("I91").Value = 2980
$ws.Range("J91").Value = 2566.6667
$ws.Range("K91").Value = 2980
$ws.Range("L91").Value = 2566.6667
$ws.Range("M91").Value = -1576
$ws.Range("N91").Value = -5374.6667
# row 136
$ws.Range("H136").Value = 1395.8654
$ws.Range("I136").Value = 1253.4615
$ws.Range("K136").Value = 3760.3845
$ws.Range("M136").Value = -1210.3845

$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 1807.64
$ws.Range("I86").Value = 1594.2632
$ws.Range("J86").Value = 2483.3333
$ws.Range("K86").Value = 1594.2632
$ws.Range("L86").Value = 2483.3333
$ws.Range("M86").Value = -471.2632000000001
$ws.Range("N86").Value = -4729.3333
# row 89
$ws.Range("H89").Value = 1807.64
$ws.Range("I89").Value = 1594.2632
$ws.Range("J89").Value = 2483.3333
$ws.Range("K89").Value = 7971.316000000001
$ws.Range("L89").Value = 12416.6665
$ws.Range("M89").Value = -2355.316000000001
$ws.Range("N89").Value = -23648.6665
# row 134
$ws.Range("H134").Value = 638088.5
$ws.Range("I134").Value = 1028586.7
$ws.Range("J134").Value = 3529
$ws.Range("K134").Value = 3085760.1
$ws.Range("L134").Value = 10587
$ws.Range("M134").Value = -3083225.1
$ws.Range("N134").Value = -15657

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 7354701
$ws.Range("I31").Value = 1104.2894
$ws.Range("J31").Value = 16669257
$ws.Range("K31").Value = 1104.2894
$ws.Range("L31").Value = 16669257
$ws.Range("M31").Value = -809.2893999999999
$ws.Range("N31").Value = -16669847
# row 34
$ws.Range("H34").Value = 7354701
$ws.Range("I34").Value = 1104.2894
$ws.Range("J34").Value = 16669257
$ws.Range("K34").Value = 1104.2894
$ws.Range("L34").Value = 16669257
$ws.Range("M34").Value = -902.2893999999999
$ws.Range("N34").Value = -16669661
# row 132
$ws.Range("H132").Value = 1159591.2
$ws.Range("I132").Value = 2240.1482
$ws.Range("J132").Value = 7409287
$ws.Range("K132").Value = 6720.444600000001
$ws.Range("L132").Value = 22227861
$ws.Range("M132").Value = -4190.444600000001
$ws.Range("N132").Value = -22232921
# row 141
$ws.Range("H141").Value = 82666.664
$ws.Range("J141").Value = 82666.664
$ws.Range("L141").Value = 82666.664
$ws.Range("N141").Value = -93026.664

$ws = $wb.Worksheets.Item("CUL")
# row 56
$ws.Range("H56").Value = 4302.857
$ws.Range("I56").Value = 4302.857
$ws.Range("K56").Value = 4302.857
$ws.Range("M56").Value = -3772.857
# row 100
$ws.Range("H100").Value = 3800
$ws.Range("J100").Value = 3800
$ws.Range("L100").Value = 11400
$ws.Range("N100").Value = -13022
# row 131
$ws.Range("H131").Value = 914.66
$ws.Range("J131").Value = 920.6804
$ws.Range("L131").Value = 2762.0412
$ws.Range("N131").Value = -12842.0412
# row 133
$ws.Range("H133").Value = 4720.654
$ws.Range("I133").Value = 2173.7
$ws.Range("J133").Value = 6312.5
$ws.Range("K133").Value = 6521.099999999999
$ws.Range("L133").Value = 18937.5
$ws.Range("M133").Value = -1461.099999999999
$ws.Range("N133").Value = -29057.5
# row 134
$ws.Range("H134").Value = 8484.647000000001
$ws.Range("I134").Value = 7885.364
$ws.Range("J134").Value = 9583.333000000001
$ws.Range("K134").Value = 23656.092
$ws.Range("L134").Value = 28749.999
$ws.Range("M134").Value = -18586.092
$ws.Range("N134").Value = -38889.999

$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 4666.575
$ws.Range("I70").Value = 4448.086
$ws.Range("J70").Value = 6196
$ws.Range("K70").Value = 4448.086
$ws.Range("L70").Value = 6196
$ws.Range("M70").Value = -4178.086
$ws.Range("N70").Value = -6736
# row 73
$ws.Range("H73").Value = 4666.575
$ws.Range("I73").Value = 4448.086
$ws.Range("J73").Value = 6196
$ws.Range("K73").Value = 4448.086
$ws.Range("L73").Value = 6196
$ws.Range("M73").Value = -3512.086
$ws.Range("N73").Value = -8068

$ws = $wb.Worksheets.Item("LTW")
# row 4
$ws.Range("H4").Value = 60006
$ws.Range("I4").Value = 30000
$ws.Range("J4").Value = 80010
$ws.Range("K4").Value = 30000
$ws.Range("L4").Value = 80010
$ws.Range("M4").Value = -29887
$ws.Range("N4").Value = -80236
# row 28
$ws.Range("H28").Value = 60006
$ws.Range("I28").Value = 30000
$ws.Range("J28").Value = 80010
$ws.Range("K28").Value = 30000
$ws.Range("L28").Value = 80010
$ws.Range("M28").Value = -29768
$ws.Range("N28").Value = -80474
# row 37
$ws.Range("H37").Value = 60006
$ws.Range("I37").Value = 30000
$ws.Range("J37").Value = 80010
$ws.Range("K37").Value = 30000
$ws.Range("L37").Value = 80010
$ws.Range("M37").Value = -29893
$ws.Range("N37").Value = -80224
# row 68
$ws.Range("H68").Value = 11220.417
$ws.Range("I68").Value = 26987.5
$ws.Range("J68").Value = 3336.875
$ws.Range("K68").Value = 26987.5
$ws.Range("L68").Value = 3336.875
$ws.Range("M68").Value = -26238.5
$ws.Range("N68").Value = -4834.875
# row 71
$ws.Range("H71").Value = 11220.417
$ws.Range("I71").Value = 26987.5
$ws.Range("J71").Value = 3336.875
$ws.Range("K71").Value = 134937.5
$ws.Range("L71").Value = 16684.375
$ws.Range("M71").Value = -131193.5
$ws.Range("N71").Value = -24172.375
# row 82
$ws.Range("H82").Value = 2147.2727
$ws.Range("I82").Value = 1892.5
$ws.Range("J82").Value = 2826.6667
$ws.Range("K82").Value = 1892.5
$ws.Range("L82").Value = 2826.6667
$ws.Range("M82").Value = -1531.5
$ws.Range("N82").Value = -3548.6667
# row 85
$ws.Range("H85").Value = 2147.2727
$ws.Range("I85").Value = 1892.5
$ws.Range("J85").Value = 2826.6667
$ws.Range("K85").Value = 1892.5
$ws.Range("L85").Value = 2826.6667
$ws.Range("M85").Value = -644.5
$ws.Range("N85").Value = -5322.6667
# row 132
$ws.Range("H132").Value = 5028.0835
$ws.Range("I132").Value = 5415.2383
$ws.Range("J132").Value = 4486.067
$ws.Range("K132").Value = 16245.7149
$ws.Range("L132").Value = 13458.201
$ws.Range("M132").Value = -13715.7149
$ws.Range("N132").Value = -18518.201
# row 136
$ws.Range("H136").Value = 1866.641
$ws.Range("I136").Value = 1325.6086
$ws.Range("J136").Value = 2644.375
$ws.Range("K136").Value = 3976.8258
$ws.Range("L136").Value = 7933.125
$ws.Range("M136").Value = -1426.8258
$ws.Range("N136").Value = -13033.125

$ws = $wb.Worksheets.Item("WVR")
# row 113
$ws.Range("H113").Value = 664.85
$ws.Range("I113").Value = 516.4666999999999
$ws.Range("J113").Value = 1110
$ws.Range("K113").Value = 1549.4001
$ws.Range("L113").Value = 3330
$ws.Range("M113").Value = 620.5999000000002
$ws.Range("N113").Value = -7670
# row 132
$ws.Range("H132").Value = 2147.3794
$ws.Range("I132").Value = 2208.878
$ws.Range("J132").Value = 1999.0588
$ws.Range("K132").Value = 6626.634
$ws.Range("L132").Value = 5997.1764
$ws.Range("M132").Value = -4096.634
$ws.Range("N132").Value = -11057.1764
